$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card10")

# Remove trailing space from the existing "Correction " header
$ws.Range("N1").Value = "Correction"

# Add the new "Serviced by " header in column O, matching the header
# formatting (bold, bordered, centered) used by the rest of row 1.
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)
$ws.Range("O1").Value = "Serviced by "
$excel.CutCopyMode = 0

# Fill in the previously-blank "Correction" cells with the same "nan"
# placeholder used throughout the rest of the table, and create the new
# (still-blank) "Serviced by " cells alongside them.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 14).Value = "nan"

    $oCell = $ws.Cells.Item($r, 15)
    $oCell.Value = "'"
    $oCell.Style = "Normal"
}
